$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A: update id numbers (1..10 -> 52501..52510)
$aValues = @(52501,52502,52503,52504,52505,52506,52507,52508,52509,52510)
for ($i = 0; $i -lt $aValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $aValues[$i]
}

# Column B: fix the misspelled name "Prretika Shetty" -> "Preetika Shetty"
$ws.Range("B3").Value = "Preetika Shetty"

# Column C: update score values
$cValues = @(80,80,97,96,93,95,92,85,91,90)
for ($i = 0; $i -lt $cValues.Length; $i++) {
    $ws.Cells.Item($i + 1, 3).Value = $cValues[$i]
}

# Update the active selection to B3
$ws.Range("B3").Select()
